# Apply updated betting odds values (FlashScore 2024-10-17 weekly games).
# Commit: "Atualizando o arquivo XLSX" -- refreshed odds for several matches.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Corinthians vs Athletico-PR
$ws.Range("Q2").Value = 1.89
$ws.Range("R2").Value = 2.01

# Row 3 - Flamengo RJ vs Fluminense
$ws.Range("G3").Value = 1.91
$ws.Range("O3").Value = 1.62
$ws.Range("P3").Value = 2.3
$ws.Range("Q3").Value = 2.88
$ws.Range("R3").Value = 1.4
$ws.Range("Z3").Value = 15
$ws.Range("AN3").Value = 3.6
$ws.Range("AO3").Value = 11
$ws.Range("AZ3").Value = 126

# Row 6 - RFS vs Auda
$ws.Range("I6").Value = 7.5
$ws.Range("K6").Value = 2.47
$ws.Range("Q6").Value = 1.53
$ws.Range("R6").Value = 2.2
$ws.Range("S6").Value = 1.26
$ws.Range("T6").Value = 3.48
$ws.Range("U6").Value = 1.82
$ws.Range("V6").Value = 1.94
$ws.Range("AB6").Value = 19
$ws.Range("AE6").Value = 15.5
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 19
$ws.Range("AI6").Value = 45
$ws.Range("AK6").Value = 150
$ws.Range("AO6").Value = 5.8
$ws.Range("AQ6").Value = 15
$ws.Range("AT6").Value = 3.35
$ws.Range("AU6").Value = 8
$ws.Range("AV6").Value = 65
$ws.Range("BA6").Value = 250
$ws.Range("BB6").Value = 450

# Row 10 - Sport Huancayo vs Grau
$ws.Range("G10").Value = 2.45
$ws.Range("I10").Value = 2.9
$ws.Range("J10").Value = 3.25
$ws.Range("L10").Value = 3.6
$ws.Range("N10").Value = 7.5
$ws.Range("U10").Value = 1.91
$ws.Range("V10").Value = 1.8
$ws.Range("W10").Value = 7
$ws.Range("X10").Value = 11
$ws.Range("Y10").Value = 10
$ws.Range("Z10").Value = 23
$ws.Range("AA10").Value = 23
$ws.Range("AH10").Value = 8
$ws.Range("AI10").Value = 13
$ws.Range("AJ10").Value = 11
$ws.Range("AK10").Value = 29
$ws.Range("AL10").Value = 26
$ws.Range("AO10").Value = 15
$ws.Range("AQ10").Value = 51
$ws.Range("AR10").Value = 81
$ws.Range("AW10").Value = 4.75
$ws.Range("AX10").Value = 17
$ws.Range("AZ10").Value = 51
$ws.Range("BA10").Value = 81

# Row 11 - Cerro Largo vs Maldonado
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 4
$ws.Range("L11").Value = 5
$ws.Range("M11").Value = 1.13
$ws.Range("N11").Value = 6
$ws.Range("W11").Value = 5
$ws.Range("Z11").Value = 17
$ws.Range("AI11").Value = 19

# Row 13 - Nacional vs Miramar
$ws.Range("G13").Value = 1.25
$ws.Range("M13").Value = 1.04
$ws.Range("N13").Value = 12
$ws.Range("AJ13").Value = 29
$ws.Range("AN13").Value = 3.2
$ws.Range("AO13").Value = 5.5
$ws.Range("AW13").Value = 11
